$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: add I0 and IF columns, matching the style of the existing header cells ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data for columns I (I0) and J (IF), rows 2-34 ---
$data = @(
    @(1,3),
    @(1,5),
    @(1,5),
    @(1,8),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,7),
    @(1,7),
    @(1,6),
    @(1,5),
    @(1,6),
    @(1,7),
    @(1,4),
    @(1,6),
    @(1,5),
    @(1,5),
    @(1,6),
    @(1,5),
    @(1,7),
    @(1,9),
    @(1,5),
    @(1,8),
    @(1,8),
    @(1,4),
    @(1,5),
    @(1,7),
    @(1,6),
    @(1,5),
    @(5,6),
    @(3,4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}

Write-Output "done"
